# Append a new data row (row 24) to the Adafruit IO data sheet,
# mirroring the existing rows' layout: Timestamp, Feed Key, Value, Latitude, Longitude, Elevation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Format the Value column as text first so the numeric-looking reading
# ("25") is stored as text, matching the rest of the sheet.
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
